$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: clear cells that must disappear entirely in the new layout ---
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()

# --- Step 2: for brand-new cells, copy number format/style from a stable donor cell first ---
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: set final values for all target cells (rows 10-23) ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5840793 - Sérgio Schneider"
$ws.Range("C15").Value = "5840793 - Sérgio Schneider"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C18").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("C19").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("C20").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("C21").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)
"
$ws.Range("C23").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)
"

# --- Step 4: row heights to match the new layout ---
$ws.Rows("13:13").RowHeight = 60
$ws.Rows("14:14").RowHeight = 60
$ws.Rows("15:15").RowHeight = 120
$ws.Rows("16:16").RowHeight = 120
$ws.Rows("17:17").AutoFit()
$ws.Rows("18:18").RowHeight = 60
$ws.Rows("19:19").RowHeight = 60
$ws.Rows("20:20").RowHeight = 60
$ws.Rows("21:21").RowHeight = 120
$ws.Rows("22:22").AutoFit()
$ws.Rows("23:23").RowHeight = 30

# --- Step 5: delete the now-extra trailing rows 24-26 ---
$ws.Rows("24:26").Delete()